$d = $word.ActiveDocument

# Change 1: expand the "aktivitásunkat (Látogatott oldalak)." sentence with
# details about blog visits/comments and comment deletion.
$d.Content.Find.Execute(
    "az aktivitásunkat (Látogatott oldalak). ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "az aktivitásunkat; blog látogatásainkat és a hozzászólásainkat. Itt lehetőségünk van hozzászólásaink törlésére is. ",
    2)

# Change 2: add "címét" (address) to the list of stored user attributes.
$d.Content.Find.Execute(
    "személyigazolvány számát, utolsó bejelentkezésének idejét, bejelentkezések számát",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "személyigazolvány számát, címét, utolsó bejelentkezésének idejét, bejelentkezések számát",
    2)
